# "added Paper of User Views"
# Adds a description for the existing "user description" row (D8) and a new
# "user views" row (C9/D9) with its description to the Work Documentation
# sheet, then leaves the selection where the user's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing "user description" entry (C8) was missing its explanation in D8.
$ws.Range("D8").Value = "Beschreibung aller in Frage kommenden Nutzer der Anwendung"

# New row: "user views" + its explanation.
$ws.Range("C9").Value = "user views"
$ws.Range("D9").Value = "Beschreibung von notwendigen Ansichten und damit verbundene Use Cases"

# Match the author's final selection after adding the new row.
[void]$ws.Range("C12").Select()
